$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to text so numeric-looking strings
# (e.g. "1.001", "30.289.32") are preserved verbatim, matching the
# source data which stores them as inline strings, not numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.289.32'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '1.870.10'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '234.99'
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').Value = '0.4697'
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('D8').Value = '0.2867'
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').Value = '0.06602'
$ws.Range('E9').Value = '  +0.95%  '
$ws.Range('D10').Value = '21.62'
$ws.Range('E10').Value = '  -1.49%  '
$ws.Range('D11').Value = '0.07961'
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').Value = '96.66'
$ws.Range('E12').Value = '  -0.66%  '
$ws.Range('D13').Value = '1.878.43'
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').Value = '0.6987'
$ws.Range('E14').Value = '  +2.25%  '
$ws.Range('D15').Value = '5.109'
$ws.Range('E15').Value = '  -1.12%  '
$ws.Range('D16').Value = '268.30'
$ws.Range('E16').Value = '  -0.87%  '
$ws.Range('D17').Value = '30.347.04'
$ws.Range('E17').Value = '  +0.42%  '
$ws.Range('D18').Value = '14.12'
$ws.Range('E18').Value = '  +4.24%  '
$ws.Range('D19').Value = '0.000007788'
$ws.Range('E19').Value = '  +5.97%  '
$ws.Range('D20').Value = '1.0000'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').Value = '2.122.08'
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').Value = '5.269'
$ws.Range('E23').Value = '  -1.17%  '
$ws.Range('D24').Value = '6.213'
$ws.Range('E24').Value = '  +0.65%  '
$ws.Range('D25').Value = '9.373'
$ws.Range('E25').Value = '  +1.58%  '
$ws.Range('D26').Value = '167.43'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').Value = '18.86'
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('D28').Value = '1.952'
$ws.Range('E28').Value = '  -0.62%  '
$ws.Range('D29').Value = '1.366'
$ws.Range('E29').Value = '  -1.25%  '
$ws.Range('D30').Value = '0.09912'
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('D31').Value = '4.335'
$ws.Range('E31').Value = '  -0.69%  '
$ws.Range('D32').Value = '1.461'
$ws.Range('E32').Value = '  -1.15%  '
$ws.Range('D33').Value = '4.053'
$ws.Range('E33').Value = '  -0.17%  '
$ws.Range('D34').Value = '0.04719'
$ws.Range('E34').Value = '  +0.22%  '
$ws.Range('D35').Value = '1.135'
$ws.Range('E35').Value = '  +0.22%  '
$ws.Range('D36').Value = '0.7025'
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').Value = '2.727'
$ws.Range('E37').Value = '  +0.62%  '
$ws.Range('D38').Value = '0.01873'
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('D39').Value = '2.796'
$ws.Range('E39').Value = '  +6.49%  '
$ws.Range('D40').Value = '6.258'
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D41').Value = '71.85'
$ws.Range('E41').Value = '  -5.00%  '
$ws.Range('D42').Value = '1.958'
$ws.Range('E42').Value = '  +0.61%  '
$ws.Range('D43').Value = '0.4179'
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('D44').Value = '0.8410'
$ws.Range('E44').Value = '  -1.13%  '
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('D46').Value = '102.77'
$ws.Range('E46').Value = '  -0.63%  '
$ws.Range('D47').Value = '7.106'
$ws.Range('E47').Value = '  -0.99%  '
$ws.Range('D48').Value = '9.173'
$ws.Range('E48').Value = '  -0.59%  '
$ws.Range('D49').Value = '917.82'
$ws.Range('E49').Value = '  -4.07%  '
$ws.Range('D50').Value = '34.57'
$ws.Range('E50').Value = '  +1.23%  '
$ws.Range('D51').Value = '0.05685'
$ws.Range('E51').Value = '  +0.61%  '
